$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in column D stay as exact text (matches source data which
# stores all Price/Volume cells as inline strings, not numbers).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = '@'
    $cell.Value = $val
}

$ws.Range('D2').Value = '68.194.46'
$ws.Range('E2').Value = '  +1.60%  '

$ws.Range('D3').Value = '2.520.19'
$ws.Range('E3').Value = '  +1.30%  '

Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.10%  '

Set-TextValue 'D5' '590.14'
$ws.Range('E5').Value = '  +1.12%  '

Set-TextValue 'D6' '177.69'
$ws.Range('E6').Value = '  +3.61%  '

$ws.Range('E7').Value = '  -0.02%  '

Set-TextValue 'D8' '0.517'
$ws.Range('E8').Value = '  +0.73%  '

$ws.Range('E9').Value = '  +3.22%  '

Set-TextValue 'D10' '0.164'
$ws.Range('E10').Value = '  -0.41%  '

$ws.Range('E11').Value = '  +2.29%  '

Set-TextValue 'D12' '4.96'
$ws.Range('E12').Value = '  +0.82%  '

$ws.Range('B13').Value = 'Avalanche'
$ws.Range('C13').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D13' '25.87'
$ws.Range('E13').Value = '  +2.08%  '

$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.876.51'
$ws.Range('E14').Value = '  -3.26%  '

$ws.Range('D15').Value = '68.013.41'
$ws.Range('E15').Value = '  +1.34%  '

$ws.Range('E16').Value = '  +1.45%  '

$ws.Range('D17').Value = '2.539.08'
$ws.Range('E17').Value = '  +3.20%  '

$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D18' '11.07'
$ws.Range('E18').Value = '  +0.40%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D19' '7.59'
$ws.Range('E19').Value = '  +2.58%  '

Set-TextValue 'D20' '353.25'
$ws.Range('E20').Value = '  +1.43%  '

Set-TextValue 'D21' '4.06'
$ws.Range('E21').Value = '  +0.66%  '

$ws.Range('E22').Value = '  +0.21%  '

Set-TextValue 'D23' '70.72'
$ws.Range('E23').Value = '  +3.24%  '

Set-TextValue 'D24' '4.32'
$ws.Range('E24').Value = '  +2.23%  '

Set-TextValue 'D25' '1.78'
$ws.Range('E25').Value = '  -0.63%  '

Set-TextValue 'D26' '9.22'
$ws.Range('E26').Value = '  -0.95%  '

$ws.Range('D27').Value = '2.636.38'
$ws.Range('E27').Value = '  +0.82%  '

$ws.Range('E28').Value = '  -0.12%  '

$ws.Range('E29').Value = '  +1.96%  '

Set-TextValue 'D30' '512.34'
$ws.Range('E30').Value = '  +0.06%  '

Set-TextValue 'D31' '7.91'
$ws.Range('E31').Value = '  +2.25%  '

$ws.Range('E32').Value = '  +4.13%  '

$ws.Range('E33').Value = '  +1.16%  '

$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('E35').Value = '  +5.39%  '

Set-TextValue 'D36' '164.31'
$ws.Range('E36').Value = '  +2.52%  '

$ws.Range('E37').Value = '  +1.10%  '

$ws.Range('E39').Value = '  +1.39%  '

$ws.Range('E41').Value = '  +3.76%  '

Set-TextValue 'D42' '4.92'
$ws.Range('E42').Value = '  +2.37%  '

$ws.Range('E43').Value = '  +0.86%  '

$ws.Range('E44').Value = '  +5.42%  '

Set-TextValue 'D45' '147.00'
$ws.Range('E45').Value = '  +3.09%  '

$ws.Range('E46').Value = '  +3.11%  '

$ws.Range('E47').Value = '  +1.53%  '

$ws.Range('E48').Value = '  +4.11%  '

Set-TextValue 'D49' '0.0746'
$ws.Range('E49').Value = '  +1.81%  '

Set-TextValue 'D50' '1.61'
$ws.Range('E50').Value = '  +2.52%  '

$ws.Range('E51').Value = '  +0.93%  '
